# Changes after Jan 23 lecture: add a new "W3" week marker column (I) plus a
# restated/duplicated marker column (J) next to the existing week/assignment
# markers in column H, shifting the week numbering down (W3 -> W7 etc. are
# unaffected text-wise; only new I/J cells are introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value  = "W3"
$ws.Range("J6").Value  = "P2"

$ws.Range("I8").Value  = "P2"
$ws.Range("J8").Value  = "W3"

$ws.Range("I9").Value  = "W4"
$ws.Range("J9").Value  = "P3"

$ws.Range("I10").Value = "W5"
$ws.Range("J10").Value = "W4"

$ws.Range("I13").Value = "P3"
$ws.Range("J13").Value = "P4"

$ws.Range("I14").Value = "W6"
$ws.Range("J14").Value = "W5"

$ws.Range("I15").Value = "W7"

$ws.Range("I16").Value = "P4"

# Update the active cell / selection to match the post-edit state.
$ws.Range("J15").Select()
